$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observations to append: date label (column A) and observed USD value (column B)
$newRows = @(
    @{ Date = "13-10-2021"; Value = 827.5599999999999 },
    @{ Date = "14-10-2021"; Value = 820.1900000000001 },
    @{ Date = "15-10-2021"; Value = 816.8099999999999 },
    @{ Date = "18-10-2021"; Value = 826.1900000000001 }
)

# Locate the last populated row in column A, then append right after it
$xlUp = [Microsoft.Office.Interop.Excel.XlDirection]::xlUp
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

$r = $lastRow + 1
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.Date
    $ws.Cells.Item($r, 2).Value = $row.Value
    $r = $r + 1
}
